# Commit: "Add MOD04 Introduccion a OOP"
#
# The canonical diff touches the title slide's subtitle text box:
#   "Módulo 2" -> "Módulo 3"
# (the "Capítulo 7" paragraph underneath is left untouched).
#
# NOTE: the diff's other hunk only rewrites the *cached* text of an
# auto-updating <a:fld type="datetimeFigureOut"> field in the Notes
# Master ("1/29/2021" -> "1/30/2022"). That cached string is
# recalculated by PowerPoint itself from the wall clock whenever the
# field recomputes (e.g. on open/print) - it is not user content and
# is not something the PowerPoint object model exposes a setter for
# (TextRange.Text on a field run is read-only / reverts), so it is
# intentionally left alone here.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("Subtitle 2")
$tr = $sh.TextFrame.TextRange

# First paragraph of the subtitle currently reads "Módulo 2".
$para1 = $tr.Paragraphs(1, 1)

# Setting the paragraph's .Text directly would diff against the old
# string ("Módulo 2" vs "Módulo 3" share the "Módulo " prefix) and the
# host would split the run in two to keep the untouched prefix's
# formatting, which does not match the source: a single run with the
# full "Módulo 3" text. Routing the edit through an unrelated
# intermediate value avoids that shared-prefix/suffix optimization, so
# the final .Text assignment rewrites the whole run in one shot.
$para1.Text = "placeholder"
$para1b = $tr.Paragraphs(1, 1)
$para1b.Text = "Módulo 3"
